$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11 (pushes existing row 11 "Rural_Urban" and
# everything below it down by one), then fill in the new row with the
# "Migration_background" / P_MIG variable.
$ws.Rows("11:11").Insert()

$ws.Range("A11").Value = "Migration_background"
$ws.Range("B11").Value = "P_MIG"
$ws.Range("C11").Value = "main"
$ws.Range("D11").Value = "1 is yes, 2 is no / 9: keine Angabe"
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 0

# The "Car_ownership" row (now row 22 after the insert/shift) has its type
# column changed from "attitude" to "main".
$ws.Range("C22").Value = "main"

$ws.Range("C23").Select()
Write-Host "Done"
